$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update want-to-go counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 109
$wsExhibit.Range("F3").Value = 946

# Sheet "全部类型" (All Types) - mirrors the same data, update identically
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 109
$wsAll.Range("F3").Value = 946
